# Updated cryptos list on Mon Oct 23 15:20:03 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Values are written as TEXT (matching the
# workbook's existing inline-string cells for these columns), using
# NumberFormat "@" + a Style reset so numeric-looking strings (like
# "219.24" or "0.999") are not silently coerced into Number cells by
# Excel's literal-value parser.
$updates = [ordered]@{
    'D2' = '30.820.76'
    'E2' = '  +3.14%  '
    'D3' = '1.678.81'
    'E3' = '  +3.17%  '
    'E4' = '  -0.29%  '
    'D5' = '219.24'
    'E5' = '  +2.30%  '
    'D6' = '0.530'
    'E6' = '  +2.32%  '
    'E7' = '  -0.31%  '
    'D8' = '29.18'
    'E8' = '  +2.55%  '
    'E9' = '  +2.40%  '
    'D10' = '0.0644'
    'E10' = '  +6.06%  '
    'D11' = '0.0905'
    'E11' = '  +0.57%  '
    'D12' = '1.919.40'
    'E12' = '  +3.17%  '
    'D13' = '1.681.81'
    'E13' = '  +3.43%  '
    'D14' = '10.18'
    'E14' = '  +10.41%  '
    'D15' = '0.603'
    'E15' = '  +7.71%  '
    'D16' = '4.04'
    'E16' = '  +5.71%  '
    'D17' = '30.798.84'
    'E17' = '  +2.98%  '
    'D18' = '65.94'
    'E18' = '  +3.21%  '
    'D19' = '243.52'
    'E19' = '  +1.34%  '
    'D20' = '0.0₃0720'
    'E20' = '  +2.85%  '
    'D21' = '0.999'
    'E21' = '  -0.19%  '
    'D22' = '4.24'
    'E22' = '  +3.07%  '
    'E23' = '  +2.28%  '
    'D24' = '2.16'
    'E24' = '  +0.22%  '
    'D25' = '159.30'
    'D26' = '15.82'
    'E26' = '  +2.59%  '
    'E27' = '  +2.38%  '
    'D28' = '6.69'
    'E28' = '  +2.43%  '
    'D29' = '0.998'
    'E29' = '  -0.36%  '
    'E30' = '  +1.52%  '
    'D31' = '1.15'
    'E31' = '  +4.16%  '
    'E32' = '  +3.09%  '
    'D33' = '1.524.88'
    'E33' = '  +7.12%  '
    'D34' = '3.30'
    'E34' = '  +4.41%  '
    'D35' = '1.76'
    'E35' = '  +5.50%  '
    'B36' = 'Aave'
    'C36' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D36' = '84.15'
    'E36' = '  +12.57%  '
    'B37' = 'TrustWalletToken'
    'C37' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D37' = '1.02'
    'E37' = '  +0.92%  '
    'D38' = '0.602'
    'E38' = '  +8.97%  '
    'D39' = '0.0179'
    'E39' = '  +5.13%  '
    'D40' = '2.66'
    'E40' = '  -3.24%  '
    'E41' = '  -0.16%  '
    'D42' = '2.03'
    'E42' = '  +2.61%  '
    'D43' = '0.838'
    'E43' = '  +1.55%  '
    'E44' = '  +0.34%  '
    'E45' = '  +2.22%  '
    'E47' = '  +3.89%  '
    'D48' = '1.811.21'
    'E48' = '  +2.41%  '
    'D49' = '50.39'
    'E49' = '  +2.78%  '
    'B50' = 'Quant'
    'C50' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D50' = '92.80'
    'E50' = '  +2.43%  '
    'B51' = 'BabyDogeCoin'
    'C51' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D51' = '0.0₆0115'
    'E51' = '  +3.44%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}

